{"js": "const replacements = [\n  [\"378\u00d79=\", \"238\u00d72=\"],\n  [\"633\u00d75=\", \"250\u00d72=\"],\n  [\"787\u00d79=\", \"680\u00d73=\"],\n  [\"207\u00d78=\", \"731\u00d72=\"],\n  [\"805\u00d73=\", \"894\u00d74=\"],\n  [\"495\u00d79=\", \"253\u00d75=\"],\n  [\"922\u00d78=\", \"236\u00d78=\"],\n  [\"136\u00d78=\", \"996\u00d75=\"],\n  [\"982\u00d79=\", \"966\u00d79=\"],\n  [\"206\u00d78=\", \"269\u00d77=\"],\n  [\"459\u00d77=\", \"526\u00d78=\"],\n  [\"271\u00d78=\", \"356\u00d76=\"],\n  [\"361\u00d73=\", \"297\u00d79=\"],\n  [\"752\u00d72=\", \"395\u00d77=\"],\n  [\"457\u00d74=\", \"687\u00d77=\"],\n  [\"871\u00d79=\", \"467\u00d76=\"],\n  [\"310\u00d78=\", \"923\u00d76=\"],\n  [\"762\u00d72=\", \"202\u00d75=\"],\n  [\"117\u00d79=\", \"792\u00d74=\"],\n  [\"455\u00d79=\", \"921\u00d72=\"],\n  [\"445\u00d75=\", \"426\u00d77=\"],\n  [\"706\u00d78=\", \"141\u00d76=\"],\n  [\"905\u00d79=\", \"337\u00d78=\"],\n  [\"491\u00d74=\", \"261\u00d77=\"],\n  [\"356\u00d74=\", \"768\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"378\u00d79=\"; New = \"238\u00d72=\" },\n    @{ Old = \"633\u00d75=\"; New = \"250\u00d72=\" },\n    @{ Old = \"787\u00d79=\"; New = \"680\u00d73=\" },\n    @{ Old = \"207\u00d78=\"; New = \"731\u00d72=\" },\n    @{ Old = \"805\u00d73=\"; New = \"894\u00d74=\" },\n    @{ Old = \"495\u00d79=\"; New = \"253\u00d75=\" },\n    @{ Old = \"922\u00d78=\"; New = \"236\u00d78=\" },\n    @{ Old = \"136\u00d78=\"; New = \"996\u00d75=\" },\n    @{ Old = \"982\u00d79=\"; New = \"966\u00d79=\" },\n    @{ Old = \"206\u00d78=\"; New = \"269\u00d77=\" },\n    @{ Old = \"459\u00d77=\"; New = \"526\u00d78=\" },\n    @{ Old = \"271\u00d78=\"; New = \"356\u00d76=\" },\n    @{ Old = \"361\u00d73=\"; New = \"297\u00d79=\" },\n    @{ Old = \"752\u00d72=\"; New = \"395\u00d77=\" },\n    @{ Old = \"457\u00d74=\"; New = \"687\u00d77=\" },\n    @{ Old = \"871\u00d79=\"; New = \"467\u00d76=\" },\n    @{ Old = \"310\u00d78=\"; New = \"923\u00d76=\" },\n    @{ Old = \"762\u00d72=\"; New = \"202\u00d75=\" },\n    @{ Old = \"117\u00d79=\"; New = \"792\u00d74=\" },\n    @{ Old = \"455\u00d79=\"; New = \"921\u00d72=\" },\n    @{ Old = \"445\u00d75=\"; New = \"426\u00d77=\" },\n    @{ Old = \"706\u00d78=\"; New = \"141\u00d76=\" },\n    @{ Old = \"905\u00d79=\"; New = \"337\u00d78=\" },\n    @{ Old = \"491\u00d74=\"; New = \"261\u00d77=\" },\n    @{ Old = \"356\u00d74=\"; New = \"768\u00d74=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute(\n        $pair.Old,     # FindText\n        $true,         # MatchCase\n        $false,        # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap (wdFindContinue)\n        $false,        # Format\n        $pair.New,     # ReplaceWith\n        2              # Replace (wdReplaceAll)\n    )\n}\n"}
